$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 10.911025
$ws.Cells.Item(2, 8).Value = 32.733075
$ws.Cells.Item(2, 9).Value = 0.03114956057965708
$ws.Cells.Item(2, 10).Value = 0.03274614149636444
$ws.Cells.Item(2, 13).Value = 12.431794
$ws.Cells.Item(2, 14).Value = 37.295382
$ws.Cells.Item(2, 15).Value = 0.6267040910788743
$ws.Cells.Item(2, 16).Value = 0.7075740515758999
$ws.Cells.Item(2, 17).Value = 135.64361512885
$ws.Cells.Item(2, 18).Value = 1220.79253615965
$ws.Cells.Item(2, 19).Value = 0.01952155705058032
$ws.Cells.Item(2, 20).Value = 0.02317032001206029
$ws.Cells.Item(3, 7).Value = 10.911025
$ws.Cells.Item(3, 8).Value = 32.733075
$ws.Cells.Item(3, 9).Value = 0.03114956057965708
$ws.Cells.Item(3, 10).Value = 0.03274614149636444
$ws.Cells.Item(3, 15).Value = 0.0264162940991436
$ws.Cells.Item(3, 16).Value = 0.0298250554119953
$ws.Cells.Item(3, 17).Value = 5.717533491358332
$ws.Cells.Item(3, 18).Value = 51.457801422225
$ws.Cells.Item(3, 19).Value = 0.0008228559533313114
$ws.Cells.Item(3, 20).Value = 0.0009766554846581082
$ws.Cells.Item(4, 7).Value = 10.911025
$ws.Cells.Item(4, 8).Value = 32.733075
$ws.Cells.Item(4, 9).Value = 0.03114956057965708
$ws.Cells.Item(4, 10).Value = 0.03274614149636444
$ws.Cells.Item(4, 13).Value = 0.03915333333333333
$ws.Cells.Item(4, 14).Value = 0.11746
$ws.Cells.Item(4, 15).Value = 0.001973774193762771
$ws.Cells.Item(4, 16).Value = 0.002228470219130754
$ws.Cells.Item(4, 17).Value = 0.4272029988333333
$ws.Cells.Item(4, 18).Value = 3.8448269895
$ws.Cells.Item(4, 19).Value = 0.00006148219881917723
$ws.Cells.Item(4, 20).Value = 0.00007297380111608995
$ws.Cells.Item(5, 7).Value = 10.911025
$ws.Cells.Item(5, 8).Value = 32.733075
$ws.Cells.Item(5, 9).Value = 0.03114956057965708
$ws.Cells.Item(5, 10).Value = 0.03274614149636444
$ws.Cells.Item(5, 13).Value = 6.8015495
$ws.Cells.Item(5, 14).Value = 13.603099
$ws.Cells.Item(5, 15).Value = 0.3428756056708687
$ws.Cells.Item(5, 16).Value = 0.2580802061075034
$ws.Cells.Item(5, 17).Value = 74.21187663323751
$ws.Cells.Item(5, 18).Value = 445.271259799425
$ws.Cells.Item(5, 19).Value = 0.01068042445013134
$ws.Cells.Item(5, 20).Value = 0.008451130946607203
$ws.Cells.Item(6, 7).Value = 10.911025
$ws.Cells.Item(6, 8).Value = 32.733075
$ws.Cells.Item(6, 9).Value = 0.03114956057965708
$ws.Cells.Item(6, 10).Value = 0.03274614149636444
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.04027333333333333
$ws.Cells.Item(6, 14).Value = 0.12082
$ws.Cells.Item(6, 15).Value = 0.002030234957350741
$ws.Cells.Item(6, 16).Value = 0.002292216685470609
$ws.Cells.Item(6, 17).Value = 0.4394233468333333
$ws.Cells.Item(6, 18).Value = 3.9548101215
$ws.Cells.Item(6, 19).Value = 0.00006324092679493441
$ws.Cells.Item(6, 20).Value = 0.00007506125192274808
$ws.Cells.Item(7, 9).Value = 0.3601508510009905
$ws.Cells.Item(7, 10).Value = 0.3786105006764219
$ws.Cells.Item(7, 13).Value = 12.431794
$ws.Cells.Item(7, 14).Value = 37.295382
$ws.Cells.Item(7, 15).Value = 0.6267040910788743
$ws.Cells.Item(7, 16).Value = 0.7075740515758999
$ws.Cells.Item(7, 17).Value = 1568.309873796748
$ws.Cells.Item(7, 18).Value = 14114.78886417073
$ws.Cells.Item(7, 19).Value = 0.2257080117278588
$ws.Cells.Item(7, 20).Value = 0.2678949659327958
$ws.Cells.Item(8, 9).Value = 0.3601508510009905
$ws.Cells.Item(8, 10).Value = 0.3786105006764219
$ws.Cells.Item(8, 15).Value = 0.0264162940991436
$ws.Cells.Item(8, 16).Value = 0.0298250554119953
$ws.Cells.Item(8, 19).Value = 0.00951385080009901
$ws.Cells.Item(8, 20).Value = 0.01129207916223757
$ws.Cells.Item(9, 9).Value = 0.3601508510009905
$ws.Cells.Item(9, 10).Value = 0.3786105006764219
$ws.Cells.Item(9, 13).Value = 0.03915333333333333
$ws.Cells.Item(9, 14).Value = 0.11746
$ws.Cells.Item(9, 15).Value = 0.001973774193762771
$ws.Cells.Item(9, 16).Value = 0.002228470219130754
$ws.Cells.Item(9, 17).Value = 4.939316019773333
$ws.Cells.Item(9, 18).Value = 44.45384417796
$ws.Cells.Item(9, 19).Value = 0.0007108564555674558
$ws.Cells.Item(9, 20).Value = 0.0008437222254075904
$ws.Cells.Item(10, 9).Value = 0.3601508510009905
$ws.Cells.Item(10, 10).Value = 0.3786105006764219
$ws.Cells.Item(10, 13).Value = 6.8015495
$ws.Cells.Item(10, 14).Value = 13.603099
$ws.Cells.Item(10, 15).Value = 0.3428756056708687
$ws.Cells.Item(10, 16).Value = 0.2580802061075034
$ws.Cells.Item(10, 17).Value = 858.0368398935291
$ws.Cells.Item(10, 18).Value = 5148.221039361174
$ws.Cells.Item(10, 19).Value = 0.1234869411698434
$ws.Cells.Item(10, 20).Value = 0.097711876049036
$ws.Cells.Item(11, 9).Value = 0.3601508510009905
$ws.Cells.Item(11, 10).Value = 0.3786105006764219
$ws.Cells.Item(11, 11).Value = 1
$ws.Cells.Item(11, 12).Value = 0.3333333333333333
$ws.Cells.Item(11, 13).Value = 0.04027333333333333
$ws.Cells.Item(11, 14).Value = 0.12082
$ws.Cells.Item(11, 15).Value = 0.002030234957350741
$ws.Cells.Item(11, 16).Value = 0.002292216685470609
$ws.Cells.Item(11, 17).Value = 5.080607538813333
$ws.Cells.Item(11, 18).Value = 45.72546784932
$ws.Cells.Item(11, 19).Value = 0.0007311908476218289
$ws.Cells.Item(11, 20).Value = 0.0008678573069448756
$ws.Cells.Item(12, 7).Value = 48.19780633333333
$ws.Cells.Item(12, 8).Value = 144.593419
$ws.Cells.Item(12, 9).Value = 0.1375984830193998
$ws.Cells.Item(12, 10).Value = 0.1446511382757993
$ws.Cells.Item(12, 13).Value = 12.431794
$ws.Cells.Item(12, 14).Value = 37.295382
$ws.Cells.Item(12, 15).Value = 0.6267040910788743
$ws.Cells.Item(12, 16).Value = 0.7075740515758999
$ws.Cells.Item(12, 17).Value = 599.1851995878953
$ws.Cells.Item(12, 18).Value = 5392.666796291058
$ws.Cells.Item(12, 19).Value = 0.08623353223450483
$ws.Cells.Item(12, 20).Value = 0.1023513919748731
$ws.Cells.Item(13, 7).Value = 48.19780633333333
$ws.Cells.Item(13, 8).Value = 144.593419
$ws.Cells.Item(13, 9).Value = 0.1375984830193998
$ws.Cells.Item(13, 10).Value = 0.1446511382757993
$ws.Cells.Item(13, 15).Value = 0.0264162940991436
$ws.Cells.Item(13, 16).Value = 0.0298250554119953
$ws.Cells.Item(13, 17).Value = 25.25634135389077
$ws.Cells.Item(13, 18).Value = 227.3070721850169
$ws.Cells.Item(13, 19).Value = 0.00363484199503648
$ws.Cells.Item(13, 20).Value = 0.00431422821448391
$ws.Cells.Item(14, 7).Value = 48.19780633333333
$ws.Cells.Item(14, 8).Value = 144.593419
$ws.Cells.Item(14, 9).Value = 0.1375984830193998
$ws.Cells.Item(14, 10).Value = 0.1446511382757993
$ws.Cells.Item(14, 13).Value = 0.03915333333333333
$ws.Cells.Item(14, 14).Value = 0.11746
$ws.Cells.Item(14, 15).Value = 0.001973774193762771
$ws.Cells.Item(14, 16).Value = 0.002228470219130754
$ws.Cells.Item(14, 17).Value = 1.887104777304444
$ws.Cells.Item(14, 18).Value = 16.98394299574
$ws.Cells.Item(14, 19).Value = 0.000271588334884596
$ws.Cells.Item(14, 20).Value = 0.0003223507538109836
$ws.Cells.Item(15, 7).Value = 48.19780633333333
$ws.Cells.Item(15, 8).Value = 144.593419
$ws.Cells.Item(15, 9).Value = 0.1375984830193998
$ws.Cells.Item(15, 10).Value = 0.1446511382757993
$ws.Cells.Item(15, 13).Value = 6.8015495
$ws.Cells.Item(15, 14).Value = 13.603099
$ws.Cells.Item(15, 15).Value = 0.3428756056708687
$ws.Cells.Item(15, 16).Value = 0.2580802061075034
$ws.Cells.Item(15, 17).Value = 327.8197655675801
$ws.Cells.Item(15, 18).Value = 1966.918593405481
$ws.Cells.Item(15, 19).Value = 0.04717916320466943
$ws.Cells.Item(15, 20).Value = 0.03733159557990327
$ws.Cells.Item(16, 7).Value = 48.19780633333333
$ws.Cells.Item(16, 8).Value = 144.593419
$ws.Cells.Item(16, 9).Value = 0.1375984830193998
$ws.Cells.Item(16, 10).Value = 0.1446511382757993
$ws.Cells.Item(16, 11).Value = 1
$ws.Cells.Item(16, 12).Value = 0.3333333333333333
$ws.Cells.Item(16, 13).Value = 0.04027333333333333
$ws.Cells.Item(16, 14).Value = 0.12082
$ws.Cells.Item(16, 15).Value = 0.002030234957350741
$ws.Cells.Item(16, 16).Value = 0.002292216685470609
$ws.Cells.Item(16, 17).Value = 1.941086320397778
$ws.Cells.Item(16, 18).Value = 17.46977688358
$ws.Cells.Item(16, 19).Value = 0.0002793572503044177
$ws.Cells.Item(16, 20).Value = 0.0003315717527281035
$ws.Cells.Item(17, 7).Value = 51.234875
$ws.Cells.Item(17, 8).Value = 102.46975
$ws.Cells.Item(17, 9).Value = 0.1462689199780642
$ws.Cells.Item(17, 10).Value = 0.102510654211286
$ws.Cells.Item(17, 13).Value = 12.431794
$ws.Cells.Item(17, 14).Value = 37.295382
$ws.Cells.Item(17, 15).Value = 0.6267040910788743
$ws.Cells.Item(17, 16).Value = 0.7075740515758999
$ws.Cells.Item(17, 17).Value = 636.9414116157501
$ws.Cells.Item(17, 18).Value = 3821.6484696945
$ws.Cells.Item(17, 19).Value = 0.09166733054794132
$ws.Cells.Item(17, 20).Value = 0.07253387892997573
$ws.Cells.Item(18, 7).Value = 51.234875
$ws.Cells.Item(18, 8).Value = 102.46975
$ws.Cells.Item(18, 9).Value = 0.1462689199780642
$ws.Cells.Item(18, 10).Value = 0.102510654211286
$ws.Cells.Item(18, 15).Value = 0.0264162940991436
$ws.Cells.Item(18, 16).Value = 0.0298250554119953
$ws.Cells.Item(18, 17).Value = 26.84780886654166
$ws.Cells.Item(18, 18).Value = 161.08685319925
$ws.Cells.Item(18, 19).Value = 0.003863882807704645
$ws.Cells.Item(18, 20).Value = 0.003057385942171495
$ws.Cells.Item(19, 7).Value = 51.234875
$ws.Cells.Item(19, 8).Value = 102.46975
$ws.Cells.Item(19, 9).Value = 0.1462689199780642
$ws.Cells.Item(19, 10).Value = 0.102510654211286
$ws.Cells.Item(19, 13).Value = 0.03915333333333333
$ws.Cells.Item(19, 14).Value = 0.11746
$ws.Cells.Item(19, 15).Value = 0.001973774193762771
$ws.Cells.Item(19, 16).Value = 0.002228470219130754
$ws.Cells.Item(19, 17).Value = 2.006016139166667
$ws.Cells.Item(19, 18).Value = 12.036096835
$ws.Cells.Item(19, 19).Value = 0.0002887018196022549
$ws.Cells.Item(19, 20).Value = 0.0002284419400534615
$ws.Cells.Item(20, 7).Value = 51.234875
$ws.Cells.Item(20, 8).Value = 102.46975
$ws.Cells.Item(20, 9).Value = 0.1462689199780642
$ws.Cells.Item(20, 10).Value = 0.102510654211286
$ws.Cells.Item(20, 13).Value = 6.8015495
$ws.Cells.Item(20, 14).Value = 13.603099
$ws.Cells.Item(20, 15).Value = 0.3428756056708687
$ws.Cells.Item(20, 16).Value = 0.2580802061075034
$ws.Cells.Item(20, 17).Value = 348.4765384388125
$ws.Cells.Item(20, 18).Value = 1393.90615375525
$ws.Cells.Item(20, 19).Value = 0.05015204452830259
$ws.Cells.Item(20, 20).Value = 0.0264559707670637
$ws.Cells.Item(21, 7).Value = 51.234875
$ws.Cells.Item(21, 8).Value = 102.46975
$ws.Cells.Item(21, 9).Value = 0.1462689199780642
$ws.Cells.Item(21, 10).Value = 0.102510654211286
$ws.Cells.Item(21, 11).Value = 1
$ws.Cells.Item(21, 12).Value = 0.3333333333333333
$ws.Cells.Item(21, 13).Value = 0.04027333333333333
$ws.Cells.Item(21, 14).Value = 0.12082
$ws.Cells.Item(21, 15).Value = 0.002030234957350741
$ws.Cells.Item(21, 16).Value = 0.002292216685470609
$ws.Cells.Item(21, 17).Value = 2.063399199166667
$ws.Cells.Item(21, 18).Value = 12.380395195
$ws.Cells.Item(21, 19).Value = 0.0002969602745134041
$ws.Cells.Item(21, 20).Value = 0.0002349766320216178
$ws.Cells.Item(22, 7).Value = 113.7817686666667
$ws.Cells.Item(22, 8).Value = 341.3453060000001
$ws.Cells.Item(22, 9).Value = 0.3248321854218885
$ws.Cells.Item(22, 10).Value = 0.3414815653401283
$ws.Cells.Item(22, 13).Value = 12.431794
$ws.Cells.Item(22, 14).Value = 37.295382
$ws.Cells.Item(22, 15).Value = 0.6267040910788743
$ws.Cells.Item(22, 16).Value = 0.7075740515758999
$ws.Cells.Item(22, 17).Value = 1414.511509019655
$ws.Cells.Item(22, 18).Value = 12730.60358117689
$ws.Cells.Item(22, 19).Value = 0.203573659517989
$ws.Cells.Item(22, 20).Value = 0.241623494726195
$ws.Cells.Item(23, 7).Value = 113.7817686666667
$ws.Cells.Item(23, 8).Value = 341.3453060000001
$ws.Cells.Item(23, 9).Value = 0.3248321854218885
$ws.Cells.Item(23, 10).Value = 0.3414815653401283
$ws.Cells.Item(23, 15).Value = 0.0264162940991436
$ws.Cells.Item(23, 16).Value = 0.0298250554119953
$ws.Cells.Item(23, 17).Value = 59.62327765335088
$ws.Cells.Item(23, 18).Value = 536.6094988801581
$ws.Cells.Item(23, 19).Value = 0.008580862542972153
$ws.Cells.Item(23, 20).Value = 0.01018470660844422
$ws.Cells.Item(24, 7).Value = 113.7817686666667
$ws.Cells.Item(24, 8).Value = 341.3453060000001
$ws.Cells.Item(24, 9).Value = 0.3248321854218885
$ws.Cells.Item(24, 10).Value = 0.3414815653401283
$ws.Cells.Item(24, 13).Value = 0.03915333333333333
$ws.Cells.Item(24, 14).Value = 0.11746
$ws.Cells.Item(24, 15).Value = 0.001973774193762771
$ws.Cells.Item(24, 16).Value = 0.002228470219130754
$ws.Cells.Item(24, 17).Value = 4.454935515862222
$ws.Cells.Item(24, 18).Value = 40.09441964276
$ws.Cells.Item(24, 19).Value = 0.0006411453848892869
$ws.Cells.Item(24, 20).Value = 0.0007609814987426287
$ws.Cells.Item(25, 7).Value = 113.7817686666667
$ws.Cells.Item(25, 8).Value = 341.3453060000001
$ws.Cells.Item(25, 9).Value = 0.3248321854218885
$ws.Cells.Item(25, 10).Value = 0.3414815653401283
$ws.Cells.Item(25, 13).Value = 6.8015495
$ws.Cells.Item(25, 14).Value = 13.603099
$ws.Cells.Item(25, 15).Value = 0.3428756056708687
$ws.Cells.Item(25, 16).Value = 0.2580802061075034
$ws.Cells.Item(25, 17).Value = 773.8923317838825
$ws.Cells.Item(25, 18).Value = 4643.353990703295
$ws.Cells.Item(25, 19).Value = 0.1113770323179219
$ws.Cells.Item(25, 20).Value = 0.08812963276489319
$ws.Cells.Item(26, 7).Value = 113.7817686666667
$ws.Cells.Item(26, 8).Value = 341.3453060000001
$ws.Cells.Item(26, 9).Value = 0.3248321854218885
$ws.Cells.Item(26, 10).Value = 0.3414815653401283
$ws.Cells.Item(26, 11).Value = 1
$ws.Cells.Item(26, 12).Value = 0.3333333333333333
$ws.Cells.Item(26, 13).Value = 0.04027333333333333
$ws.Cells.Item(26, 14).Value = 0.12082
$ws.Cells.Item(26, 15).Value = 0.002030234957350741
$ws.Cells.Item(26, 16).Value = 0.002292216685470609
$ws.Cells.Item(26, 17).Value = 4.582371096768889
$ws.Cells.Item(26, 18).Value = 41.24133987092001
$ws.Cells.Item(26, 19).Value = 0.0006594856581161558
$ws.Cells.Item(26, 20).Value = 0.0007827497418532642
